$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to stay text so numeric-looking values
# (e.g. "1.26", "0.997") are not silently reinterpreted as numbers by
# Excel's smart cell entry -- the source data is textual (e.g. "58.960.21",
# "  -6.04%  ") and must round-trip as strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '58.960.21'
$ws.Range("E2").Value = '  -6.04%  '
$ws.Range("D3").Value = '2.447.20'
$ws.Range("E3").Value = '  -8.51%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '538.92'
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("D6").Value = '145.38'
$ws.Range("E6").Value = '  -7.37%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").Value = '0.571'
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").Value = '2.460.86'
$ws.Range("E9").Value = '  -8.02%  '
$ws.Range("D10").Value = '0.0991'
$ws.Range("E10").Value = '  -5.99%  '
$ws.Range("E11").Value = '  -1.75%  '
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '0.349'
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("D14").Value = '2.881.12'
$ws.Range("E14").Value = '  -8.53%  '
$ws.Range("D15").Value = '23.93'
$ws.Range("E15").Value = '  -8.88%  '
$ws.Range("D16").Value = '58.850.13'
$ws.Range("E16").Value = '  -6.12%  '
$ws.Range("D17").Value = '0.0000138'
$ws.Range("E17").Value = '  -5.17%  '
$ws.Range("D18").Value = '2.502.67'
$ws.Range("E18").Value = '  -6.62%  '
$ws.Range("D19").Value = '11.14'
$ws.Range("E19").Value = '  -5.06%  '
$ws.Range("D20").Value = '4.34'
$ws.Range("E20").Value = '  -5.38%  '
$ws.Range("D21").Value = '324.63'
$ws.Range("E21").Value = '  -5.52%  '
$ws.Range("D22").Value = '0.966'
$ws.Range("E22").Value = '  -3.34%  '
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  -7.36%  '
$ws.Range("D24").Value = '60.37'
$ws.Range("E24").Value = '  -4.38%  '
$ws.Range("D25").Value = '0.451'
$ws.Range("E25").Value = '  -11.26%  '
$ws.Range("D26").Value = '0.160'
$ws.Range("E26").Value = '  -5.18%  '
$ws.Range("D27").Value = '0.975'
$ws.Range("E27").Value = '  -2.34%  '
$ws.Range("D28").Value = '7.71'
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("D29").Value = '1.82'
$ws.Range("E29").Value = '  -5.77%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.26'
$ws.Range("E30").Value = '  -8.71%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0770'
$ws.Range("E31").Value = '  -8.97%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '6.63'
$ws.Range("E32").Value = '  -8.35%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Value = '155.74'
$ws.Range("E34").Value = '  -4.81%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.39'
$ws.Range("E35").Value = '  -4.42%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '18.53'
$ws.Range("E36").Value = '  -4.28%  '
$ws.Range("D37").Value = '4.44'
$ws.Range("E37").Value = '  -8.35%  '
$ws.Range("D38").Value = '1.69'
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").Value = '5.83'
$ws.Range("E39").Value = '  -4.80%  '
$ws.Range("D40").Value = '314.70'
$ws.Range("E40").Value = '  -7.70%  '
$ws.Range("D41").Value = '36.20'
$ws.Range("E41").Value = '  -5.41%  '
$ws.Range("D42").Value = '0.833'
$ws.Range("E42").Value = '  -10.65%  '
$ws.Range("D43").Value = '3.70'
$ws.Range("E43").Value = '  -6.57%  '
$ws.Range("D44").Value = '0.995'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").Value = '10.73'
$ws.Range("E45").Value = '  -2.57%  '
$ws.Range("D46").Value = '0.590'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").Value = '0.0935'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("D48").Value = '0.0526'
$ws.Range("E48").Value = '  -4.64%  '
$ws.Range("D49").Value = '123.01'
$ws.Range("E49").Value = '  -3.67%  '
$ws.Range("D50").Value = '0.0230'
$ws.Range("E50").Value = '  -4.29%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '18.89'
$ws.Range("E51").Value = '  -8.63%  '
